$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"
$wsMeta.Range("B8").Value = "2024-06-19T17:47:42+02:00"
$wsMeta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Elements sheet updates ---
$wsElem = $wb.Worksheets.Item("Elements")

# classCode row (row 12): Min 1 -> 0 and Base Min 1 -> 0.
# Copy/PasteSpecial(values) from an existing "0" text cell keeps the
# cell as a text-typed shared string (matching the existing column
# formatting) instead of Excel auto-converting the literal "0" into a
# numeric cell.
$wsElem.Range("F3").Copy()
$wsElem.Range("F12").PasteSpecial(-4163)
$wsElem.Range("AG3").Copy()
$wsElem.Range("AG12").PasteSpecial(-4163)

# Binding Value Set URLs for classCode/moodCode
$wsElem.Range("Z12").Value = "http://hl7.org/cda/stds/core/ValueSet/CDAActClassObservation"
$wsElem.Range("Z13").Value = "http://hl7.org/cda/stds/core/ValueSet/CDAActMood"

$wsElem.Application.CutCopyMode = $false
